$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet currently has a blank/implicit row at 27 (a gap row, no explicit
# XML row) between the "50/23-24" entry (row 26) and the "2" group (rows
# 28-29). We need to insert a brand-new data row just above that "2" group,
# which pushes rows 28-31 down to 29-32, while row 27 becomes the home of
# the new entry.
#
# Inserting directly at row 28 (rather than row 27) means the row being
# "pushed down from above" is the blank/implicit row 27, so there is no
# explicit formatting on it to blend into the freshly inserted row - this
# keeps the style table clean (matches the source file, which shows no
# styles.xml changes at all).
$ws.Rows.Item(28).Insert()

# Populate the new row 27 with the "51/23-24" entry.
$ws.Cells.Item(27, 2).Value2 = 45184
$ws.Cells.Item(27, 3).Value2 = "51/23-24"
$ws.Cells.Item(27, 4).Value2 = "Namrata Rubber Product Pvt Ltd"
$ws.Cells.Item(27, 5).Value2 = 26491
$ws.Cells.Item(27, 6).Formula = "=F26+E27"

# Carry over the formatting used by the rest of this group (row 26) onto the
# new row 27, cell by cell, so the existing style entries get reused instead
# of new ones being created.
$ws.Cells.Item(26, 1).Copy()
$ws.Cells.Item(27, 1).PasteSpecial(-4122)

$ws.Cells.Item(26, 2).Copy()
$ws.Cells.Item(27, 2).PasteSpecial(-4122)

$ws.Cells.Item(26, 3).Copy()
$ws.Cells.Item(27, 3).PasteSpecial(-4122)

$ws.Cells.Item(26, 4).Copy()
$ws.Cells.Item(27, 4).PasteSpecial(-4122)

$ws.Cells.Item(26, 5).Copy()
$ws.Cells.Item(27, 5).PasteSpecial(-4122)

# F27 is the new "running total" cell for this subgroup, which carries the
# bold/bordered total style also used lower down in this same table (e.g.
# F29, the subtotal cell for the next group).
$ws.Cells.Item(29, 6).Copy()
$ws.Cells.Item(27, 6).PasteSpecial(-4122)

# Row 26's total cell is no longer the last row before a gap, so it loses
# the bold "subtotal" styling and picks up the plain running-total look used
# throughout the rest of the table (e.g. F25).
$ws.Cells.Item(25, 6).Copy()
$ws.Cells.Item(26, 6).PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the saved selection/cursor position recorded in the workbook.
$ws.Range("C37").Select()
